$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append below the existing data (15, 16, 17 marzo update).
$newRows = @(
    @(44301, 4, 44, 244.7435754811436),
    @(44302, 8, 44, 244.7435754811436),
    @(44303, 5, 39, 216.9318055401046)
)

$startRow = 227

# Copy the date-column formatting (border/bold/center/date numberformat)
# from the last existing data row onto the new date cells before filling
# in the values.
$ws.Range("A226").Copy() | Out-Null
$lastRow = $startRow + $newRows.Length - 1
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
    $ws.Cells.Item($r, 4).Value = $newRows[$i][3]
}
